$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date column G, rows 2 and 4
$overview.Range("G2").Value = "2016-08-17 20:14:25"
$overview.Range("G4").Value = "2016-08-17 20:14:25"

# zh-cn sheet: Priority column E (ht -> mt), rows 2 and 4
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E4").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime column H, rows 2 and 4
$zhcn.Range("H2").Value = "2016-08-17 20:14:20"
$zhcn.Range("H4").Value = "2016-08-17 20:14:20"

# zh-cn sheet: Correspond Handback DateTime column K, rows 2 and 4
$zhcn.Range("K2").Value = "2016-08-17 20:14:44"
$zhcn.Range("K4").Value = "2016-08-17 20:14:44"

# de-de sheet: Correspond Handoff Datetime column H, rows 2 and 4
$dede.Range("H2").Value = "2016-08-17 20:14:25"
$dede.Range("H4").Value = "2016-08-17 20:14:25"

# de-de sheet: Correspond Handback DateTime column K, rows 2 and 4
$dede.Range("K2").Value = "2016-08-17 20:14:51"
$dede.Range("K4").Value = "2016-08-17 20:14:51"
